$wb = $excel.ActiveWorkbook

# The second sheet ("21_Properties of Tangent to Cir") is renamed to a
# shorter title ("21_Prop of Tangent to Circle"). The sheet's own A1 label
# cell carries the same text as the tab name, so it needs the same update.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "21_Prop of Tangent to Circle"
$ws2.Range("A1").Value = "21_Prop of Tangent to Circle"

# The workbook's active/selected tab moves from the first sheet to this
# (renamed) second sheet, with its selection resting on the default A1 cell.
$ws2.Select()
$ws2.Range("A1").Select()
